$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New match data for rows 10-15 (Cagliari matches) that previously had
# only home/away team names filled in (columns B, C). This adds the
# remaining stats: xG_home (D), xG_away (E), goals_home (F), goals_away (G).
# All of these columns store their values as text (shared strings), matching
# the existing rows 2-9, so values are entered with a leading apostrophe to
# force text entry, then the style is reset to Normal to avoid leaving a
# "quote prefix" cell format applied (matching the original look of rows 2-9).
$data = @(
    @(10, "1.44548",  "2.7383",   "2", "2"),
    @(11, "3.02671",  "1.67924",  "1", "1"),
    @(12, "1.61081",  "2.74594",  "1", "3"),
    @(13, "0.870374", "0.265316", "0", "0"),
    @(14, "0.628497", "1.43641",  "1", "1"),
    @(15, "3.36165",  "1.19646",  "3", "2")
)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 4; $col -le 7; $col++) {
        $value = $row[$col - 3]
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    }
}
